$wb = $excel.ActiveWorkbook

# ALERTS sheet - add row 9
$ws = $wb.Worksheets.Item("ALERTS")
$ws.Cells.Item(9, 1).Value = "'2026-01-28"
$ws.Cells.Item(9, 2).Value = "15:02:42"
$ws.Cells.Item(9, 3).Value = "15:00"
$ws.Cells.Item(9, 4).Value = "Bathroom"
$ws.Cells.Item(9, 5).Value = "WARNING"
$ws.Cells.Item(9, 6).Value = "Bathroom Humidity > 90.0% for 20s with NO MOTION. Alerting."

# PIR sheet - add rows 175-187
$ws = $wb.Worksheets.Item("PIR")
$pirData = @(
    @("'2026-01-28", "15:02:08", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:10", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:16", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:16", "15:00", "Bathroom", "Motion Detected", "Active"),
    @("'2026-01-28", "15:02:24", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:28", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:33", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:38", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:44", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:48", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:53", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:02:59", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("'2026-01-28", "15:03:04", "15:00", "Bathroom", "No Motion", "Inactive")
)
$startRow = 175
for ($i = 0; $i -lt $pirData.Count; $i++) {
    $rowData = $pirData[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowData[$j]
    }
}

# Humidity sheet - add rows 169-182
$ws = $wb.Worksheets.Item("Humidity")
$humData = @(
    @("'2026-01-28", "15:02:07", "15:00", "Bathroom", "'88.3%", "Active"),
    @("'2026-01-28", "15:02:08", "15:00", "Bathroom", "'87.3%", "Active"),
    @("'2026-01-28", "15:02:11", "15:00", "Bathroom", "'88.3%", "Active"),
    @("'2026-01-28", "15:02:15", "15:00", "Bathroom", "'88.3%", "Active"),
    @("'2026-01-28", "15:02:19", "15:00", "Bathroom", "'99.9%", "Active"),
    @("'2026-01-28", "15:02:23", "15:00", "Bathroom", "'99.9%", "Active"),
    @("'2026-01-28", "15:02:31", "15:00", "Bathroom", "'99.9%", "Active"),
    @("'2026-01-28", "15:02:35", "15:00", "Bathroom", "'99.9%", "Active"),
    @("'2026-01-28", "15:02:43", "15:00", "Bathroom", "'99.9%", "Active"),
    @("'2026-01-28", "15:02:47", "15:00", "Bathroom", "'98.0%", "Active"),
    @("'2026-01-28", "15:02:51", "15:00", "Bathroom", "'95.4%", "Active"),
    @("'2026-01-28", "15:02:55", "15:00", "Bathroom", "'93.5%", "Active"),
    @("'2026-01-28", "15:02:59", "15:00", "Bathroom", "'91.4%", "Active"),
    @("'2026-01-28", "15:03:03", "15:00", "Bathroom", "'91.5%", "Active")
)
$startRow = 169
for ($i = 0; $i -lt $humData.Count; $i++) {
    $rowData = $humData[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowData[$j]
    }
}

# Temperature sheet - add rows 169-182
$ws = $wb.Worksheets.Item("Temperature")
$tempData = @(
    @("'2026-01-28", "15:02:07", "15:00", "Bathroom", "22.9C", "Active"),
    @("'2026-01-28", "15:02:09", "15:00", "Bathroom", "22.9C", "Active"),
    @("'2026-01-28", "15:02:11", "15:00", "Bathroom", "22.9C", "Active"),
    @("'2026-01-28", "15:02:15", "15:00", "Bathroom", "22.9C", "Active"),
    @("'2026-01-28", "15:02:19", "15:00", "Bathroom", "23.1C", "Active"),
    @("'2026-01-28", "15:02:23", "15:00", "Bathroom", "23.1C", "Active"),
    @("'2026-01-28", "15:02:31", "15:00", "Bathroom", "23.1C", "Active"),
    @("'2026-01-28", "15:02:35", "15:00", "Bathroom", "23.1C", "Active"),
    @("'2026-01-28", "15:02:43", "15:00", "Bathroom", "23.0C", "Active"),
    @("'2026-01-28", "15:02:47", "15:00", "Bathroom", "23.1C", "Active"),
    @("'2026-01-28", "15:02:52", "15:00", "Bathroom", "23.0C", "Active"),
    @("'2026-01-28", "15:02:56", "15:00", "Bathroom", "23.0C", "Active"),
    @("'2026-01-28", "15:03:00", "15:00", "Bathroom", "23.0C", "Active"),
    @("'2026-01-28", "15:03:04", "15:00", "Bathroom", "23.0C", "Active")
)
$startRow = 169
for ($i = 0; $i -lt $tempData.Count; $i++) {
    $rowData = $tempData[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowData[$j]
    }
}

# mmWave sheet - add rows 7-8
$ws = $wb.Worksheets.Item("mmWave")
$mmData = @(
    @("'2026-01-28", "15:02:09", "15:00", "Living Room", "Presence Detected", "Active"),
    @("'2026-01-28", "15:02:54", "15:00", "Living Room", "Presence Detected", "Active")
)
$startRow = 7
for ($i = 0; $i -lt $mmData.Count; $i++) {
    $rowData = $mmData[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowData[$j]
    }
}
